$d = $word.ActiveDocument

# New white / tiny Calibri formatting shared by every line in the (new) list.
$WHITE = 16777215   # RGB(255,255,255) == 0x00FFFFFF, BGR-packed OLE_COLOR
$FONTNAME = "Calibri"
$FONTSIZE = 4        # half-points 8 -> 4 pt

function Set-HiddenCitationFormat($range) {
    $range.Font.Name = $FONTNAME
    $range.Font.NameAscii = $FONTNAME
    $range.Font.Color = $WHITE
    $range.Font.Size = $FONTSIZE
}

function Set-ParagraphText($paragraphIndex, $text) {
    $p = $d.Paragraphs($paragraphIndex)
    $r = $p.Range
    $r.Text = $text
    # Exclude the trailing paragraph mark so the rPr lands only on the run,
    # not on the paragraph's pPr.
    $body = $d.Range($r.Start, $r.End - 1)
    Set-HiddenCitationFormat $body
}

# 1) First paragraph: title -> citation line (leading space kept).
Set-ParagraphText 1 " Willard and Spackman's occupational therapy"

# 2) Second paragraph: title -> citation line (leading space kept).
Set-ParagraphText 2 " Pedretti's Occupational Therapy-E-Book: Practice Skills for Physical Dysfunction"

# 3) Six additional citation lines appended as new paragraphs.
$newLines = @(
    " Concepts of occupational therapy",
    " Conceptual foundations of occupational therapy practice",
    " Occupational therapy for children",
    " A model of human occupation: Theory and application",
    "Adult norms for the Box and Block Test of manual dexterity",
    "Relative contributions of neural mechanisms versus muscle mechanics in promoting finger extension deficits following stroke"
)

foreach ($line in $newLines) {
    $last = $d.Paragraphs($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter()
    $newIndex = $d.Paragraphs.Count
    Set-ParagraphText $newIndex $line
}
